$d = $word.ActiveDocument

# 1. "feedback that comes from" -> "feedback that comes with"
$d.Content.Find.Execute("feedback that comes from", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "feedback that comes with", 2) | Out-Null

# 2. "one directional persuasive communication." -> "one one way persuasive communication."
$d.Content.Find.Execute("one directional persuasive communication.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "one one way persuasive communication.", 2) | Out-Null

# 3. " generating maximum revenue as creative thinking" -> " generating maximum revenue in minimum time as creative thinking"
$d.Content.Find.Execute("generating maximum revenue as creative thinking", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "generating maximum revenue in minimum time as creative thinking", 2) | Out-Null

# 4. Move the hidden "_GoBack" bookmark from after "...can damage the" to right after "one one way"
$oldBm = $d.Bookmarks("_GoBack")
if ($oldBm.Exists) {
    $oldBm.Delete()
}

$r = $d.Content
$r.Find.Execute("one one way", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($r.Find.Found) {
    $target = $d.Range($r.End, $r.End)
    $d.Bookmarks.Add("_GoBack", $target) | Out-Null
}
